$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "49.480.73"
$ws.Cells.Item(2, 5).Value = "  -1.05%  "

$ws.Cells.Item(3, 4).Value = "2.627.95"
$ws.Cells.Item(3, 5).Value = "  -0.72%  "

$ws.Cells.Item(4, 5).Value = "  +0.09%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "112.14"
$ws.Cells.Item(5, 5).Value = "  +1.07%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "323.84"
$ws.Cells.Item(6, 5).Value = "  -1.30%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.526"
$ws.Cells.Item(7, 5).Value = "  -1.04%  "

$ws.Cells.Item(8, 5).Value = "  +0.03%  "

$ws.Cells.Item(9, 5).Value = "  -3.14%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "39.58"
$ws.Cells.Item(10, 5).Value = "  -2.97%  "

$ws.Cells.Item(11, 5).Value = "  -4.23%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0810"
$ws.Cells.Item(12, 5).Value = "  -1.42%  "

$ws.Cells.Item(13, 5).Value = "  +1.24%  "

$ws.Cells.Item(14, 5).Value = "  -0.29%  "

$ws.Cells.Item(15, 4).Value = "2.990.95"
$ws.Cells.Item(15, 5).Value = "  -1.93%  "

$ws.Cells.Item(16, 4).Value = "2.624.22"
$ws.Cells.Item(16, 5).Value = "  +0.57%  "

$ws.Cells.Item(17, 5).Value = "  -3.05%  "

$ws.Cells.Item(18, 4).Value = "49.404.08"
$ws.Cells.Item(18, 5).Value = "  -1.05%  "

$ws.Cells.Item(19, 5).Value = "  -3.69%  "

$ws.Cells.Item(20, 5).Value = "  -4.18%  "

$ws.Cells.Item(21, 5).Value = "  -2.46%  "

$ws.Cells.Item(22, 5).Value = "  -2.24%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "268.99"
$ws.Cells.Item(23, 5).Value = "  -4.91%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "68.96"
$ws.Cells.Item(24, 5).Value = "  -5.77%  "

$ws.Cells.Item(25, 5).Value = "  -2.66%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "26.20"
$ws.Cells.Item(26, 5).Value = "  -2.40%  "

$ws.Cells.Item(27, 5).Value = "  +0.08%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "10.27"
$ws.Cells.Item(28, 5).Value = "  +2.90%  "

$ws.Cells.Item(29, 5).Value = "  -0.73%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.137"
$ws.Cells.Item(30, 5).Value = "  -4.84%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "34.64"
$ws.Cells.Item(31, 5).Value = "  -5.78%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "49.41"
$ws.Cells.Item(32, 5).Value = "  -0.73%  "

$ws.Cells.Item(33, 5).Value = "  +0.78%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0812"
$ws.Cells.Item(34, 5).Value = "  +1.81%  "

$ws.Cells.Item(35, 5).Value = "  -0.11%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "18.80"
$ws.Cells.Item(36, 5).Value = "  -3.80%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "4.88"
$ws.Cells.Item(37, 5).Value = "  +2.48%  "

$ws.Cells.Item(38, 5).Value = "  -1.68%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.10"

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "128.41"
$ws.Cells.Item(40, 5).Value = "  +2.77%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.110"
$ws.Cells.Item(41, 5).Value = "  -2.06%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "22.17"
$ws.Cells.Item(42, 5).Value = "  -1.43%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0325"
$ws.Cells.Item(43, 5).Value = "  +3.58%  "

$ws.Cells.Item(44, 5).Value = "  -3.80%  "

$ws.Cells.Item(45, 4).Value = "2.049.08"
$ws.Cells.Item(45, 5).Value = "  -0.94%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.20"
$ws.Cells.Item(46, 5).Value = "  -5.23%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.09"
$ws.Cells.Item(47, 5).Value = "  +4.57%  "

$ws.Cells.Item(48, 5).Value = "  -5.39%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "8.89"
$ws.Cells.Item(49, 5).Value = "  -2.16%  "

$ws.Cells.Item(50, 2).Value = "THORChain"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "5.19"
$ws.Cells.Item(50, 5).Value = "  -4.03%  "

$ws.Cells.Item(51, 2).Value = "MultiversX"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "58.67"
$ws.Cells.Item(51, 5).Value = "  +0.89%  "
